$d = $word.ActiveDocument

# The "Quota" table-schema bullet currently reads:
#   Quota (QuotaID, CompanyID, WarehouseID, ProductID, QuotaAmount , Month, IsHidden)
# A "Year" column is being added to the quota table (alongside the existing
# "Month"), so the bullet needs to become:
#   Quota (QuotaID, CompanyID, WarehouseID, ProductID, QuotaAmount , Month, Year, IsHidden)
$old = "QuotaAmount , Month, IsHidden)"
$new = "QuotaAmount , Month, Year, IsHidden)"

$found = $d.Content.Find.Execute(
    $old,   # FindText
    $true,  # MatchCase
    $false, # MatchWholeWord
    $false, # MatchWildcards
    $false, # MatchSoundsLike
    $false, # MatchAllWordForms
    $true,  # Forward
    1,      # Wrap (wdFindContinue)
    $false, # Format
    $new,   # ReplaceWith
    1)      # Replace (wdReplaceOne -- touch only this single occurrence)

if (-not $found) {
    throw "Could not find the Quota schema text to update"
}
